$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 8 ("Otel Metrics") - fill in the (currently empty) content placeholder
# ---------------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$tr8 = $s8.Shapes.Item(2).TextFrame.TextRange
$tr8.Text = "Counters`vCPU/Memory usage on the box`vRequests per second"
$tr8.Paragraphs(2,1).IndentLevel = 2
$tr8.Paragraphs(3,1).IndentLevel = 2

# ---------------------------------------------------------------------------
# Slide 9 ("Otel Traces") - fill in the (currently empty) content placeholder
# ---------------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$sp9 = $s9.Shapes.Item(2)
$tr9 = $sp9.TextFrame.TextRange

$lines9 = @(
    "Information for a request",
    "Can span multiple applications, all based on the Trace Id",
    "Traces have child Spans",
    "Spans have child Attributes",
    "Attributes are key/value pairs of data",
    "Spans have child Events",
    "Events are timestamped pieces of data in a span, can have their own attributes",
    "Example:",
    "Trace is created when a request comes into the system",
    "A span is made to wrap around call to database",
    "Span stores time it took to run the database query",
    "Another span is created when making call out to 3rd party service",
    "Span stores time it took to run the HTTP request and the response HTTP status code, plus error if needed"
)
$tr9.Text = [string]::Join("`v", $lines9)

$levels9 = @(1,2,1,2,3,2,3,1,2,2,3,2,3)
for ($i = 1; $i -le $levels9.Length; $i++) {
    $tr9.Paragraphs($i,1).IndentLevel = $levels9[$i-1]
}

# Split out the "rd" superscript in paragraph 12 ("...call out to 3rd party service")
$para12 = $tr9.Paragraphs(12,1)
$rdStart = $para12.Start + "Another span is created when making call out to 3".Length
$rdRange = $tr9.Characters($rdStart, 2)
$rdRange.Font.Superscript = $true

# Turn on "shrink text on overflow" for this placeholder
$sp9.TextFrame.AutoSize = 2

# ---------------------------------------------------------------------------
# Add 4 new slides at the end (positions 10-13), all using the
# "Title and Content" layout (the same layout used by the rest of the deck).
# ---------------------------------------------------------------------------

# --- Slide 10: "Otel Log Records" ---
$s10 = $p.Slides.Add(10, 2)
$trT10 = $s10.Shapes.Item(1).TextFrame.TextRange
$trT10.Text = "Otel"
[void]$trT10.InsertAfter(" Log Records")

$tr10 = $s10.Shapes.Item(2).TextFrame.TextRange
$lines10 = @(
    "Standalone Logs",
    "A Span with an Event, no link to other Spans",
    "Embedded Logs",
    "An Event in a Span"
)
$tr10.Text = [string]::Join("`v", $lines10)
$tr10.Paragraphs(2,1).IndentLevel = 2
$tr10.Paragraphs(4,1).IndentLevel = 2

# --- Slide 11: "OTel with C#" ---
$s11 = $p.Slides.Add(11, 2)
$s11.Shapes.Item(1).TextFrame.TextRange.Text = "OTel with C#"

$tr11 = $s11.Shapes.Item(2).TextFrame.TextRange
$lines11 = @(
    "2 APIs",
    "Because of course there are",
    "Built in to .NET one",
    "OpenTelemetry"
)
$tr11.Text = [string]::Join("`v", $lines11)
$tr11.Paragraphs(2,1).IndentLevel = 2
[void]$tr11.InsertAfter(" ")
[void]$tr11.InsertAfter("community provided one")

# --- Slide 12: blank placeholder slide ---
$s12 = $p.Slides.Add(12, 2)

# --- Slide 13: blank placeholder slide ---
$s13 = $p.Slides.Add(13, 2)
